# Adds a new data row (row 5) to Tabelle1 demonstrating the
# "ExecResult.Insights" test case: if-cond match count + analyzed
# datarow count, represented here as a "Text" / "notmodified" / "warning"
# sample row appended below the existing test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Text"
$ws.Range("B5").Value = "notmodified"
$ws.Range("C5").Value = "warning"

# Move the active selection, matching the saved workbook view state.
$ws.Range("H6").Select() | Out-Null
